$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(118, 8).Value = 595
$ws.Cells.Item(118, 9).Value = 595
$ws.Cells.Item(118, 11).Value = 1785
$ws.Cells.Item(118, 13).Value = -128

$ws.Cells.Item(129, 8).Value = 873.0755
$ws.Cells.Item(129, 9).Value = 719.2
$ws.Cells.Item(129, 10).Value = 889.1042
$ws.Cells.Item(129, 11).Value = 2157.6
$ws.Cells.Item(129, 12).Value = 2667.3126
$ws.Cells.Item(129, 13).Value = 2842.4
$ws.Cells.Item(129, 14).Value = -12667.3126

$ws.Cells.Item(132, 8).Value = 818.9
$ws.Cells.Item(132, 9).Value = 755.54346
$ws.Cells.Item(132, 11).Value = 2266.63038
$ws.Cells.Item(132, 13).Value = 263.3696199999999

$ws.Cells.Item(137, 8).Value = 1850.7931
$ws.Cells.Item(137, 9).Value = 1310.3572
$ws.Cells.Item(137, 10).Value = 2355.2
$ws.Cells.Item(137, 11).Value = 3931.0716
$ws.Cells.Item(137, 12).Value = 7065.599999999999
$ws.Cells.Item(137, 13).Value = -1381.0716
$ws.Cells.Item(137, 14).Value = -12165.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 19999.25
$ws.Cells.Item(24, 10).Value = 19999.25
$ws.Cells.Item(24, 12).Value = 19999.25
$ws.Cells.Item(24, 14).Value = -20747.25

$ws.Cells.Item(32, 8).Value = 2413.6396
$ws.Cells.Item(32, 9).Value = 1681.7106
$ws.Cells.Item(32, 11).Value = 1681.7106
$ws.Cells.Item(32, 13).Value = -1394.7106

$ws.Cells.Item(46, 8).Value = 11541
$ws.Cells.Item(46, 9).Value = 9663
$ws.Cells.Item(46, 10).Value = 12793
$ws.Cells.Item(46, 11).Value = 9663
$ws.Cells.Item(46, 12).Value = 12793
$ws.Cells.Item(46, 13).Value = -9344
$ws.Cells.Item(46, 14).Value = -13431

$ws.Cells.Item(74, 8).Value = 1246.4667
$ws.Cells.Item(74, 9).Value = 833.619
$ws.Cells.Item(74, 11).Value = 833.619
$ws.Cells.Item(74, 13).Value = 40.38099999999997

$ws.Cells.Item(77, 8).Value = 1246.4667
$ws.Cells.Item(77, 9).Value = 833.619
$ws.Cells.Item(77, 11).Value = 4168.095
$ws.Cells.Item(77, 13).Value = 199.9049999999997

$ws.Cells.Item(100, 8).Value = 19999.25
$ws.Cells.Item(100, 10).Value = 19999.25
$ws.Cells.Item(100, 12).Value = 19999.25
$ws.Cells.Item(100, 14).Value = -22163.25

$ws.Cells.Item(109, 8).Value = 61250.332
$ws.Cells.Item(109, 10).Value = 61250.332
$ws.Cells.Item(109, 12).Value = 61250.332
$ws.Cells.Item(109, 14).Value = -64024.332

$ws.Cells.Item(114, 8).Value = 9000
$ws.Cells.Item(114, 10).Value = 9000
$ws.Cells.Item(114, 12).Value = 9000
$ws.Cells.Item(114, 14).Value = -17678

$ws.Cells.Item(139, 8).Value = 51926.25
$ws.Cells.Item(139, 10).Value = 51926.25
$ws.Cells.Item(139, 12).Value = 51926.25
$ws.Cells.Item(139, 14).Value = -62206.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 60000
$ws.Cells.Item(76, 10).Value = 60000
$ws.Cells.Item(76, 12).Value = 60000
$ws.Cells.Item(76, 14).Value = -60630

$ws.Cells.Item(79, 8).Value = 60000
$ws.Cells.Item(79, 10).Value = 60000
$ws.Cells.Item(79, 12).Value = 60000
$ws.Cells.Item(79, 14).Value = -62184

$ws.Cells.Item(81, 8).Value = 18074.8
$ws.Cells.Item(81, 10).Value = 18074.8
$ws.Cells.Item(81, 12).Value = 18074.8
$ws.Cells.Item(81, 14).Value = -20196.8

$ws.Cells.Item(84, 8).Value = 18074.8
$ws.Cells.Item(84, 10).Value = 18074.8
$ws.Cells.Item(84, 12).Value = 54224.39999999999
$ws.Cells.Item(84, 14).Value = -64832.39999999999

$ws.Cells.Item(86, 8).Value = 89490.52
$ws.Cells.Item(86, 9).Value = 3308.4546
$ws.Cells.Item(86, 11).Value = 3308.4546
$ws.Cells.Item(86, 13).Value = -2185.4546

$ws.Cells.Item(89, 8).Value = 89490.52
$ws.Cells.Item(89, 9).Value = 3308.4546
$ws.Cells.Item(89, 11).Value = 16542.273
$ws.Cells.Item(89, 13).Value = -10926.273

$ws.Cells.Item(134, 8).Value = 7829.12
$ws.Cells.Item(134, 9).Value = 9043.632
$ws.Cells.Item(134, 11).Value = 27130.896
$ws.Cells.Item(134, 13).Value = -24595.896

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1892290.8
$ws.Cells.Item(58, 9).Value = 2289673
$ws.Cells.Item(58, 11).Value = 2289673
$ws.Cells.Item(58, 13).Value = -2289470

$ws.Cells.Item(74, 8).Value = 31249.5
$ws.Cells.Item(74, 10).Value = 31249.5
$ws.Cells.Item(74, 12).Value = 31249.5
$ws.Cells.Item(74, 14).Value = -32997.5

$ws.Cells.Item(77, 8).Value = 31249.5
$ws.Cells.Item(77, 10).Value = 31249.5
$ws.Cells.Item(77, 12).Value = 93748.5
$ws.Cells.Item(77, 14).Value = -102484.5

$ws.Cells.Item(86, 8).Value = 62501412
$ws.Cells.Item(86, 9).Value = 76924104
$ws.Cells.Item(86, 11).Value = 76924104
$ws.Cells.Item(86, 13).Value = -76922981

$ws.Cells.Item(89, 8).Value = 62501412
$ws.Cells.Item(89, 9).Value = 76924104
$ws.Cells.Item(89, 11).Value = 384620520
$ws.Cells.Item(89, 13).Value = -384614904

$ws.Cells.Item(96, 8).Value = 32500
$ws.Cells.Item(96, 10).Value = 32500
$ws.Cells.Item(96, 12).Value = 32500
$ws.Cells.Item(96, 14).Value = -37992

$ws.Cells.Item(132, 8).Value = 2080.3547
$ws.Cells.Item(132, 9).Value = 1312.6666
$ws.Cells.Item(132, 11).Value = 3937.9998
$ws.Cells.Item(132, 13).Value = -1407.9998

$ws.Cells.Item(134, 8).Value = 1331.6571
$ws.Cells.Item(134, 9).Value = 1361.8387
$ws.Cells.Item(134, 11).Value = 4085.5161
$ws.Cells.Item(134, 13).Value = -1550.5161

$ws.Cells.Item(136, 8).Value = 1892290.8
$ws.Cells.Item(136, 9).Value = 2289673
$ws.Cells.Item(136, 11).Value = 6869019
$ws.Cells.Item(136, 13).Value = -6866469

$ws.Cells.Item(141, 8).Value = 71948
$ws.Cells.Item(141, 10).Value = 71948
$ws.Cells.Item(141, 12).Value = 71948
$ws.Cells.Item(141, 14).Value = -82308

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 1000
$ws.Cells.Item(46, 9).Value = 1500
$ws.Cells.Item(46, 11).Value = 4500
$ws.Cells.Item(46, 13).Value = -4409

$ws.Cells.Item(131, 8).Value = 9800.225
$ws.Cells.Item(131, 10).Value = 10228.941
$ws.Cells.Item(131, 12).Value = 30686.823
$ws.Cells.Item(131, 14).Value = -40766.823

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5333.3335
$ws.Cells.Item(70, 9).Value = 7000
$ws.Cells.Item(70, 11).Value = 7000
$ws.Cells.Item(70, 13).Value = -6730

$ws.Cells.Item(73, 8).Value = 5333.3335
$ws.Cells.Item(73, 9).Value = 7000
$ws.Cells.Item(73, 11).Value = 7000
$ws.Cells.Item(73, 13).Value = -6064

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()

$ws.Cells.Item(92, 8).Value = 23299.666
$ws.Cells.Item(92, 10).Value = 23299.666
$ws.Cells.Item(92, 12).Value = 23299.666
$ws.Cells.Item(92, 14).Value = -27043.666

$ws.Cells.Item(132, 8).Value = 1835567.5
$ws.Cells.Item(132, 9).Value = 2568041.2
$ws.Cells.Item(132, 10).Value = 4383
$ws.Cells.Item(132, 11).Value = 7704123.600000001
$ws.Cells.Item(132, 12).Value = 13149
$ws.Cells.Item(132, 13).Value = -7701593.600000001
$ws.Cells.Item(132, 14).Value = -18209

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2339.6956
$ws.Cells.Item(7, 9).Value = 1610.1904
$ws.Cells.Item(7, 11).Value = 1610.1904
$ws.Cells.Item(7, 13).Value = -1498.1904

$ws.Cells.Item(16, 8).Value = 3025.353
$ws.Cells.Item(16, 9).Value = 2927.1667
$ws.Cells.Item(16, 11).Value = 2927.1667
$ws.Cells.Item(16, 13).Value = -2757.1667

$ws.Cells.Item(55, 8).Value = 347.4516
$ws.Cells.Item(55, 9).Value = 282.04
$ws.Cells.Item(55, 11).Value = 282.04
$ws.Cells.Item(55, 13).Value = -109.04

$ws.Cells.Item(82, 8).Value = 951.1
$ws.Cells.Item(82, 9).Value = 929.7143
$ws.Cells.Item(82, 10).Value = 1001
$ws.Cells.Item(82, 11).Value = 929.7143
$ws.Cells.Item(82, 12).Value = 1001
$ws.Cells.Item(82, 13).Value = -568.7143
$ws.Cells.Item(82, 14).Value = -1723

$ws.Cells.Item(85, 8).Value = 951.1
$ws.Cells.Item(85, 9).Value = 929.7143
$ws.Cells.Item(85, 10).Value = 1001
$ws.Cells.Item(85, 11).Value = 929.7143
$ws.Cells.Item(85, 12).Value = 1001
$ws.Cells.Item(85, 13).Value = 318.2857
$ws.Cells.Item(85, 14).Value = -3497

$ws.Cells.Item(126, 8).Value = 2339.6956
$ws.Cells.Item(126, 9).Value = 1610.1904
$ws.Cells.Item(126, 11).Value = 4830.5712
$ws.Cells.Item(126, 13).Value = -2360.5712

$ws.Cells.Item(127, 8).Value = 50000
$ws.Cells.Item(127, 10).Value = 50000
$ws.Cells.Item(127, 12).Value = 50000
$ws.Cells.Item(127, 14).Value = -59920

$ws.Cells.Item(132, 8).Value = 1593.579
$ws.Cells.Item(132, 9).Value = 1286.2858
$ws.Cells.Item(132, 10).Value = 2454
$ws.Cells.Item(132, 11).Value = 3858.8574
$ws.Cells.Item(132, 12).Value = 7362
$ws.Cells.Item(132, 13).Value = -1328.8574
$ws.Cells.Item(132, 14).Value = -12422

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 130962.164
$ws.Cells.Item(122, 9).Value = 156854.8
$ws.Cells.Item(122, 11).Value = 470564.4
$ws.Cells.Item(122, 13).Value = -468114.4

$ws.Cells.Item(132, 8).Value = 1972.5
$ws.Cells.Item(132, 9).Value = 1367.2667
$ws.Cells.Item(132, 11).Value = 4101.800099999999
$ws.Cells.Item(132, 13).Value = -1571.800099999999
